$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the Price (column D) value looks like a pure number and must be
# forced to text formatting so Excel keeps it as a string (matching the source data)
# instead of silently converting it to a numeric cell.
$numericPriceRows = @(5, 6, 11, 14, 15, 19, 20, 21, 22, 23, 25, 28, 29, 30, 31, 33, 35, 36, 38, 39, 40, 41, 43, 44, 50, 51)
foreach ($r in $numericPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '65.625.65'
$ws.Cells.Item(2, 5).Value = '  +0.16%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.564.27'
$ws.Cells.Item(3, 5).Value = '  +3.60%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '605.58'
$ws.Cells.Item(5, 5).Value = '  +2.25%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '141.60'
$ws.Cells.Item(6, 5).Value = '  +3.88%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.564.65'
$ws.Cells.Item(7, 5).Value = '  +3.64%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.16%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.19%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +3.07%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '7.07'
$ws.Cells.Item(11, 5).Value = '  -3.80%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +4.82%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '4.167.14'
$ws.Cells.Item(13, 5).Value = '  +3.79%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '0.0000189'
$ws.Cells.Item(14, 5).Value = '  +4.24%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '27.39'
$ws.Cells.Item(15, 5).Value = '  +3.15%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.561.81'
$ws.Cells.Item(16, 5).Value = '  +3.98%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +1.61%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '65.580.77'
$ws.Cells.Item(18, 5).Value = '  +0.20%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '10.36'
$ws.Cells.Item(19, 5).Value = '  +4.39%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '5.95'
$ws.Cells.Item(20, 5).Value = '  +2.22%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '14.39'
$ws.Cells.Item(21, 5).Value = '  +5.35%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '396.81'
$ws.Cells.Item(22, 5).Value = '  +0.73%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '0.575'
$ws.Cells.Item(23, 5).Value = '  +5.16%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '3.705.61'
$ws.Cells.Item(24, 5).Value = '  +3.48%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '74.37'
$ws.Cells.Item(25, 5).Value = '  +1.46%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.03%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +12.42%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '7.95'
$ws.Cells.Item(28, 5).Value = '  +10.16%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'PancakeSwap'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(29, 4).Value = '2.32'
$ws.Cells.Item(29, 5).Value = '  +2.92%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  +0.05%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '8.40'
$ws.Cells.Item(31, 5).Value = '  +2.63%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '3.573.27'
$ws.Cells.Item(32, 5).Value = '  +3.78%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '0.149'
$ws.Cells.Item(33, 5).Value = '  +1.45%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +0.03%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '23.92'
$ws.Cells.Item(35, 5).Value = '  +3.99%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '1.31'
$ws.Cells.Item(36, 5).Value = '  +9.58%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +2.25%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Monero'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(38, 4).Value = '169.92'
$ws.Cells.Item(38, 5).Value = '  -1.04%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).Value = '1.56'
$ws.Cells.Item(39, 5).Value = '  +4.11%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '5.08'
$ws.Cells.Item(40, 5).Value = '  +5.58%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '0.0816'
$ws.Cells.Item(41, 5).Value = '  +6.07%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +1.47%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '26.56'
$ws.Cells.Item(43, 5).Value = '  +16.46%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '43.05'
$ws.Cells.Item(44, 5).Value = '  -1.17%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.03%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +1.03%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +10.21%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +5.82%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '2.489.64'
$ws.Cells.Item(49, 5).Value = '  +13.08%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '6.86'
$ws.Cells.Item(50, 5).Value = '  +4.54%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '2.38'
$ws.Cells.Item(51, 5).Value = '  +20.26%  '
